$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 184, shifting existing rows 184-236 down to 185-237.
$ws.Rows("184").Insert()

# Populate the new row 184 with the new weekly record.
$ws.Range("A184").Value = 4
$ws.Range("B184").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C184").Value = "Los Lagos"
$ws.Range("D184").Value = 44642
$ws.Range("E184").Value = 10
$ws.Range("F184").Value = 100112003
$ws.Range("G184").Value = "Ajo"
$ws.Range("H184").Value = "Chino"
$ws.Range("I184").Value = "Primera"
$ws.Range("J184").Value = 220
$ws.Range("K184").Value = 21000
$ws.Range("L184").Value = 21000
$ws.Range("M184").Value = 21000
$ws.Range("N184").Value = "$/caja 10 kilos"
$ws.Range("O184").Value = "China"
$ws.Range("P184").Value = 2100
$ws.Range("Q184").Value = 10
$ws.Range("R184").Value = "Hortaliza"
